# Updates crypto Price (D) and Volume(1h) (E) columns per the scraper's latest run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.874.84"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "3.679.56"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("E6").Value = "  +11.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "658.15"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("E9").Value = "  +3.84%  "

$ws.Range("E10").Value = "  +0.02%  "

$ws.Range("D11").Value = "3.677.91"
$ws.Range("E11").Value = "  +2.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.70"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.44%  "

$ws.Range("D15").Value = "4.365.17"
$ws.Range("E15").Value = "  +2.43%  "

$ws.Range("E16").Value = "  +4.51%  "

$ws.Range("D17").Value = "96.687.88"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +14.86%  "

$ws.Range("D19").Value = "3.663.31"
$ws.Range("E19").Value = "  +2.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.531"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "532.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.42%  "

$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.17"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.41%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("E28").Value = "  +3.89%  "

$ws.Range("E29").Value = "  +5.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.39"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.05"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +16.81%  "

$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "664.46"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.73"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.31%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  +5.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.91"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.162"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.95%  "

$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +10.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.965"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.03"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +17.92%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0462"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.50%  "

$ws.Range("E47").Value = "  +6.19%  "

$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("E49").Value = "  +7.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.62"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.70"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.39%  "
